$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.008.91"
$ws.Range("E2").Value = "'  -0.96%  "
$ws.Range("D3").Value = "'3.419.88"
$ws.Range("E3").Value = "'  -0.75%  "
$ws.Range("E4").Value = "'  +0.18%  "
$ws.Range("D5").Value = "'410.58"
$ws.Range("E5").Value = "'  +0.61%  "
$ws.Range("D6").Value = "'129.59"
$ws.Range("E6").Value = "'  -3.30%  "
$ws.Range("E7").Value = "'  +8.33%  "
$ws.Range("E8").Value = "'  -0.06%  "
$ws.Range("E9").Value = "'  +7.67%  "
$ws.Range("E10").Value = "'  +14.86%  "
$ws.Range("D11").Value = "'42.83"
$ws.Range("E11").Value = "'  +0.83%  "
$ws.Range("D12").Value = "'0.0000217"
$ws.Range("E12").Value = "'  +63.10%  "
$ws.Range("D13").Value = "'9.13"
$ws.Range("E13").Value = "'  +7.72%  "
$ws.Range("E14").Value = "'  -0.38%  "
$ws.Range("D15").Value = "'3.954.77"
$ws.Range("E15").Value = "'  -0.75%  "
$ws.Range("D16").Value = "'21.28"
$ws.Range("E16").Value = "'  +6.55%  "
$ws.Range("D17").Value = "'3.404.68"
$ws.Range("E17").Value = "'  +1.85%  "
$ws.Range("E18").Value = "'  +5.81%  "
$ws.Range("E19").Value = "'  +5.90%  "
$ws.Range("D20").Value = "'61.989.00"
$ws.Range("E20").Value = "'  -0.85%  "
$ws.Range("D21").Value = "'443.05"
$ws.Range("E21").Value = "'  +40.56%  "
$ws.Range("D22").Value = "'91.15"
$ws.Range("E22").Value = "'  +8.19%  "
$ws.Range("E23").Value = "'  -0.62%  "
$ws.Range("D24").Value = "'13.13"
$ws.Range("E24").Value = "'  +1.22%  "
$ws.Range("D25").Value = "'3.26"
$ws.Range("E25").Value = "'  +2.90%  "
$ws.Range("D26").Value = "'33.61"
$ws.Range("E26").Value = "'  +12.59%  "
$ws.Range("D27").Value = "'8.86"
$ws.Range("E27").Value = "'  +6.53%  "
$ws.Range("E28").Value = "'  +0.61%  "
$ws.Range("D29").Value = "'7.61"
$ws.Range("E29").Value = "'  +0.26%  "
$ws.Range("E30").Value = "'  +0.43%  "
$ws.Range("D31").Value = "'12.03"
$ws.Range("E31").Value = "'  +5.52%  "
$ws.Range("E32").Value = "'  -0.41%  "
$ws.Range("E33").Value = "'  -2.87%  "
$ws.Range("D34").Value = "'42.94"
$ws.Range("E34").Value = "'  +1.38%  "
$ws.Range("E35").Value = "'  -0.13%  "
$ws.Range("D36").Value = "'0.0504"
$ws.Range("E36").Value = "'  +3.57%  "
$ws.Range("D37").Value = "'53.72"
$ws.Range("E37").Value = "'  +4.14%  "
$ws.Range("E38").Value = "'  +0.12%  "
$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "'  -1.13%  "
$ws.Range("E40").Value = "'  +7.79%  "
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "'  -1.54%  "
$ws.Range("D42").Value = "'0.315"
$ws.Range("E42").Value = "'  -1.07%  "
$ws.Range("D43").Value = "'141.42"
$ws.Range("E43").Value = "'  +3.08%  "
$ws.Range("D44").Value = "'4.23"
$ws.Range("E44").Value = "'  +4.62%  "
$ws.Range("E45").Value = "'  -0.84%  "
$ws.Range("D47").Value = "'16.67"
$ws.Range("E47").Value = "'  -1.29%  "
$ws.Range("D48").Value = "'22.22"
$ws.Range("E48").Value = "'  +3.78%  "
$ws.Range("D49").Value = "'3.767.94"
$ws.Range("D50").Value = "'2.109.35"
$ws.Range("E50").Value = "'  -1.06%  "
$ws.Range("D51").Value = "'105.33"
$ws.Range("E51").Value = "'  +25.31%  "
